# Generate Report for Handback
# Update the timestamp cells that record when the handoff/handback
# xliff files were generated/processed for the 188242da-... file.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Column G = "Latest HO Xliff Generate Date" for row 2 (188242da-...)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-18 21:06:48"

# --- zh-cn sheet ---
# H2 = Correspond Handoff Datetime, K2 = Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-18 21:06:43"
$wsZhCn.Range("K2").Value = "2016-08-18 21:06:59"

# --- de-de sheet ---
# H2 = Correspond Handoff Datetime (shares text with Overview!G2's new value),
# K2 = Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-18 21:06:48"
$wsDeDe.Range("K2").Value = "2016-08-18 21:07:12"
